# Update stats for 2026-02 (row 27 in Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B27").Value = 6559
$ws.Range("D27").Value = 6117686
$ws.Range("E27").Value = 932.7162677237384
$ws.Range("F27").Value = 10.23529411764705
$ws.Range("H27").Value = 25.47794516549766
